$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Model Accuracy" ---
$ws1 = $wb.Worksheets.Item("Model Accuracy")

# Replicate the existing bold/border style (currently on A2) down through A3:A6
# before we overwrite values, so the new rows keep the same look-and-feel.
$ws1.Range("A2").Copy($ws1.Range("A3:A6"))

$ws1.Range("A2").Value = "TOTALENERGIES SE"
$ws1.Range("B2").Value = 64.36430317848411

$ws1.Range("A3").Value = "FMC CORP"
$ws1.Range("B3").Value = 75.67237163814181

$ws1.Range("A4").Value = "BP PLC"
$ws1.Range("B4").Value = 93.7041564792176

$ws1.Range("A5").Value = "STORA ENSO"
$ws1.Range("B5").Value = 91.99266503667481

$ws1.Range("A6").Value = "BHP GROUP"
$ws1.Range("B6").Value = 93.27628361858191

# --- Sheet 2: "Signal Market Correlation" ---
$ws2 = $wb.Worksheets.Item("Signal Market Correlation")

# Drop the old "Pearson Correlation (need to stationarize)"/D-column layout;
# the sheet goes from B/C/D data columns down to just B/C.
$ws2.Range("D1").EntireColumn.Delete()

$ws2.Range("B1").Value = "Pearson Correlation"
$ws2.Range("C1").Value = "P-value"

# Replicate the existing bold/border style (currently on A2) down through A3:A6.
$ws2.Range("A2").Copy($ws2.Range("A3:A6"))

$ws2.Range("A2").Value = "TOTALENERGIES SE"
$ws2.Range("B2").Value = 0.01566256551086196
$ws2.Range("C2").Value = 0.600871449251472

$ws2.Range("A3").Value = "FMC CORP"
$ws2.Range("B3").Value = 0.01770014237732665
$ws2.Range("C3").Value = 0.5110361495306047

$ws2.Range("A4").Value = "BP PLC"
$ws2.Range("B4").Value = 0.009638066301792383
$ws2.Range("C4").Value = 0.7001561347054477

$ws2.Range("A5").Value = "STORA ENSO"
$ws2.Range("B5").Value = -0.01530076064468157
$ws2.Range("C5").Value = 0.5439945818950872

$ws2.Range("A6").Value = "BHP GROUP"
$ws2.Range("B6").Value = -0.01422839565825462
$ws2.Range("C6").Value = 0.5719727335739155

Write-Host "daily model results updated"
